$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column A so the existing columns shift right by one.
$ws.Columns.Item(1).Insert()

# Header row (row 1)
$ws.Cells.Item(1, 1).Value = "TabName"
$ws.Cells.Item(1, 2).Value = "query"
$ws.Cells.Item(1, 3).Value = "StatQuery"
$ws.Cells.Item(1, 4).Value = "dbExcel"
$ws.Cells.Item(1, 5).Value = "WebExcel"

# Data row (row 2)
$ws.Cells.Item(2, 1).Value = "CasesTab"

$casesQuery = "MATCH (ct:clinical_trial)<--(a:arm)<--(c:case)`n    WHERE c.ethnicity IN ['NOT_HISPANIC_OR_LATINO'] `nWITH DISTINCT c, a, ct`nRETURN `n    COALESCE(c.case_id, '') AS ``Case ID``,`n    COALESCE(ct.clinical_trial_designation, '') AS ``Trial Code``,`n    COALESCE(a.arm_id, '') AS ``Arm``,`n    COALESCE(a.arm_drug, '') AS ``Arm Treatment``,`n    COALESCE(c.disease, '') AS ``Diagnosis``,`n    COALESCE(c.gender, '') AS ``Gender``,`n    COALESCE(c.race, '') AS ``Race``,`n    COALESCE(c.ethnicity, '') AS ``Ethnicity``"
$ws.Cells.Item(2, 2).Value = $casesQuery

$statQuery = "MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)`n    WHERE WHERE c.ethnicity IN ['NOT_HISPANIC_OR_LATINO'] `nOPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)`nRETURN `n    COUNT(DISTINCT f) AS number_of_files,`n    COUNT(DISTINCT c.case_id) AS number_of_cases,`n    COUNT(DISTINCT ct.clinical_trial_designation) AS number_of_trials"
$ws.Cells.Item(2, 3).Value = $statQuery

$ws.Cells.Item(2, 4).Value = "TC02_Trials_Filter_Ethnicity-NotHispLatino_Neo4jData.xlsx"
$ws.Cells.Item(2, 5).Value = "TC02_Trials_Filter_Ethnicity-NotHispLatino_WebData.xlsx"

# Apply wrap-text to B2 and C2, matching the long-query cells.
$ws.Cells.Item(2, 2).WrapText = $true
$ws.Cells.Item(2, 3).WrapText = $true

# Row height for row 2
$ws.Rows.Item(2).RowHeight = 174

# Column width for the newly inserted column A (existing columns keep their widths).
$ws.Columns.Item(1).ColumnWidth = 8

# Selection matches the diff (C2 selected)
$ws.Range("C2").Select()
